$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.303422808647156
$ws.Range("B1").Value = 2.699458837509155
$ws.Range("C1").Value = 3.112922430038452
$ws.Range("D1").Value = 1.555411100387573
$ws.Range("E1").Value = 1.124618172645569
